$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Dnajb11"
$ws.Cells.Item(2,3).Value = "Prtg"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 11.25725833333333
$ws.Cells.Item(2,8).Value = 33.771775
$ws.Cells.Item(2,9).Value = 0.1213189509762199
$ws.Cells.Item(2,10).Value = 0.1213189509762199
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 0.130725
$ws.Cells.Item(2,14).Value = 0.392175
$ws.Cells.Item(2,15).Value = 0.1184385768344502
$ws.Cells.Item(2,16).Value = 0.1184385768344502
$ws.Cells.Item(2,17).Value = 1.471605095625
$ws.Cells.Item(2,18).Value = 13.244445860625
$ws.Cells.Item(2,19).Value = 0.01436884389667193
$ws.Cells.Item(2,20).Value = 0.01436884389667193

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Dnajb11"
$ws.Cells.Item(3,3).Value = "Prtg"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 11.25725833333333
$ws.Cells.Item(3,8).Value = 33.771775
$ws.Cells.Item(3,9).Value = 0.1213189509762199
$ws.Cells.Item(3,10).Value = 0.1213189509762199
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 0.7274516666666667
$ws.Cells.Item(3,14).Value = 2.182355
$ws.Cells.Item(3,15).Value = 0.6590808193983468
$ws.Cells.Item(3,16).Value = 0.659080819398347
$ws.Cells.Item(3,17).Value = 8.189111336680556
$ws.Cells.Item(3,18).Value = 73.702002030125
$ws.Cells.Item(3,19).Value = 0.07995899361795492
$ws.Cells.Item(3,20).Value = 0.0799589936179549

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Dnajb11"
$ws.Cells.Item(4,3).Value = "Prtg"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 11.25725833333333
$ws.Cells.Item(4,8).Value = 33.771775
$ws.Cells.Item(4,9).Value = 0.1213189509762199
$ws.Cells.Item(4,10).Value = 0.1213189509762199
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 0.24556
$ws.Cells.Item(4,14).Value = 0.73668
$ws.Cells.Item(4,15).Value = 0.2224806037672029
$ws.Cells.Item(4,16).Value = 0.2224806037672029
$ws.Cells.Item(4,17).Value = 2.764332356333333
$ws.Cells.Item(4,18).Value = 24.878991207
$ws.Cells.Item(4,19).Value = 0.0269911134615931
$ws.Cells.Item(4,20).Value = 0.0269911134615931

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Dnajb11"
$ws.Cells.Item(5,3).Value = "Prtg"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 18.645164
$ws.Cells.Item(5,8).Value = 55.935492
$ws.Cells.Item(5,9).Value = 0.2009380677142005
$ws.Cells.Item(5,10).Value = 0.2009380677142004
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 0.6666666666666666
$ws.Cells.Item(5,13).Value = 0.130725
$ws.Cells.Item(5,14).Value = 0.392175
$ws.Cells.Item(5,15).Value = 0.1184385768344502
$ws.Cells.Item(5,16).Value = 0.1184385768344502
$ws.Cells.Item(5,17).Value = 2.4373890639
$ws.Cells.Item(5,18).Value = 21.9365015751
$ws.Cells.Item(5,19).Value = 0.0237988187719343
$ws.Cells.Item(5,20).Value = 0.0237988187719343

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Dnajb11"
$ws.Cells.Item(6,3).Value = "Prtg"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 18.645164
$ws.Cells.Item(6,8).Value = 55.935492
$ws.Cells.Item(6,9).Value = 0.2009380677142005
$ws.Cells.Item(6,10).Value = 0.2009380677142004
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 0.7274516666666667
$ws.Cells.Item(6,14).Value = 2.182355
$ws.Cells.Item(6,15).Value = 0.6590808193983468
$ws.Cells.Item(6,16).Value = 0.659080819398347
$ws.Cells.Item(6,17).Value = 13.56345562707333
$ws.Cells.Item(6,18).Value = 122.07110064366
$ws.Cells.Item(6,19).Value = 0.1324344263173957
$ws.Cells.Item(6,20).Value = 0.1324344263173957

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Dnajb11"
$ws.Cells.Item(7,3).Value = "Prtg"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 18.645164
$ws.Cells.Item(7,8).Value = 55.935492
$ws.Cells.Item(7,9).Value = 0.2009380677142005
$ws.Cells.Item(7,10).Value = 0.2009380677142004
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 0.24556
$ws.Cells.Item(7,14).Value = 0.73668
$ws.Cells.Item(7,15).Value = 0.2224806037672029
$ws.Cells.Item(7,16).Value = 0.2224806037672029
$ws.Cells.Item(7,17).Value = 4.578506471839999
$ws.Cells.Item(7,18).Value = 41.20655824656
$ws.Cells.Item(7,19).Value = 0.04470482262487042
$ws.Cells.Item(7,20).Value = 0.04470482262487042

# Row 8
$ws.Cells.Item(8,1).Value = "M1"
$ws.Cells.Item(8,2).Value = "Dnajb11"
$ws.Cells.Item(8,3).Value = "Prtg"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 32.95867399999999
$ws.Cells.Item(8,8).Value = 98.87602199999999
$ws.Cells.Item(8,9).Value = 0.3551941011611514
$ws.Cells.Item(8,10).Value = 0.3551941011611514
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 0.6666666666666666
$ws.Cells.Item(8,13).Value = 0.130725
$ws.Cells.Item(8,14).Value = 0.392175
$ws.Cells.Item(8,15).Value = 0.1184385768344502
$ws.Cells.Item(8,16).Value = 0.1184385768344502
$ws.Cells.Item(8,17).Value = 4.308522658649999
$ws.Cells.Item(8,18).Value = 38.77670392784999
$ws.Cells.Item(8,19).Value = 0.04206868384151852
$ws.Cells.Item(8,20).Value = 0.04206868384151852

# Row 9
$ws.Cells.Item(9,1).Value = "M1"
$ws.Cells.Item(9,2).Value = "Dnajb11"
$ws.Cells.Item(9,3).Value = "Prtg"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 32.95867399999999
$ws.Cells.Item(9,8).Value = 98.87602199999999
$ws.Cells.Item(9,9).Value = 0.3551941011611514
$ws.Cells.Item(9,10).Value = 0.3551941011611514
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 0.7274516666666667
$ws.Cells.Item(9,14).Value = 2.182355
$ws.Cells.Item(9,15).Value = 0.6590808193983468
$ws.Cells.Item(9,16).Value = 0.659080819398347
$ws.Cells.Item(9,17).Value = 23.97584233242333
$ws.Cells.Item(9,18).Value = 215.78258099181
$ws.Cells.Item(9,19).Value = 0.234101619238751
$ws.Cells.Item(9,20).Value = 0.234101619238751

# Row 10
$ws.Cells.Item(10,1).Value = "M1"
$ws.Cells.Item(10,2).Value = "Dnajb11"
$ws.Cells.Item(10,3).Value = "Prtg"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 32.95867399999999
$ws.Cells.Item(10,8).Value = 98.87602199999999
$ws.Cells.Item(10,9).Value = 0.3551941011611514
$ws.Cells.Item(10,10).Value = 0.3551941011611514
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 0.24556
$ws.Cells.Item(10,14).Value = 0.73668
$ws.Cells.Item(10,15).Value = 0.2224806037672029
$ws.Cells.Item(10,16).Value = 0.2224806037672029
$ws.Cells.Item(10,17).Value = 8.09333198744
$ws.Cells.Item(10,18).Value = 72.83998788695999
$ws.Cells.Item(10,19).Value = 0.0790237980808819
$ws.Cells.Item(10,20).Value = 0.07902379808088192

# Row 11
$ws.Cells.Item(11,1).Value = "M2"
$ws.Cells.Item(11,2).Value = "Dnajb11"
$ws.Cells.Item(11,3).Value = "Prtg"
$ws.Cells.Item(11,4).Value = "ECs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 23.97572
$ws.Cells.Item(11,8).Value = 71.92716
$ws.Cells.Item(11,9).Value = 0.2583852225089954
$ws.Cells.Item(11,10).Value = 0.2583852225089954
$ws.Cells.Item(11,11).Value = 2
$ws.Cells.Item(11,12).Value = 0.6666666666666666
$ws.Cells.Item(11,13).Value = 0.130725
$ws.Cells.Item(11,14).Value = 0.392175
$ws.Cells.Item(11,15).Value = 0.1184385768344502
$ws.Cells.Item(11,16).Value = 0.1184385768344502
$ws.Cells.Item(11,17).Value = 3.134225997
$ws.Cells.Item(11,18).Value = 28.208033973
$ws.Cells.Item(11,19).Value = 0.03060277802901817
$ws.Cells.Item(11,20).Value = 0.03060277802901817

# Row 12
$ws.Cells.Item(12,1).Value = "M2"
$ws.Cells.Item(12,2).Value = "Dnajb11"
$ws.Cells.Item(12,3).Value = "Prtg"
$ws.Cells.Item(12,4).Value = "FAPs"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 23.97572
$ws.Cells.Item(12,8).Value = 71.92716
$ws.Cells.Item(12,9).Value = 0.2583852225089954
$ws.Cells.Item(12,10).Value = 0.2583852225089954
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 0.7274516666666667
$ws.Cells.Item(12,14).Value = 2.182355
$ws.Cells.Item(12,15).Value = 0.6590808193983468
$ws.Cells.Item(12,16).Value = 0.659080819398347
$ws.Cells.Item(12,17).Value = 17.44117747353333
$ws.Cells.Item(12,18).Value = 156.9705972618
$ws.Cells.Item(12,19).Value = 0.1702967441716529
$ws.Cells.Item(12,20).Value = 0.1702967441716529

# Row 13
$ws.Cells.Item(13,1).Value = "M2"
$ws.Cells.Item(13,2).Value = "Dnajb11"
$ws.Cells.Item(13,3).Value = "Prtg"
$ws.Cells.Item(13,4).Value = "sCs"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 23.97572
$ws.Cells.Item(13,8).Value = 71.92716
$ws.Cells.Item(13,9).Value = 0.2583852225089954
$ws.Cells.Item(13,10).Value = 0.2583852225089954
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 0.24556
$ws.Cells.Item(13,14).Value = 0.73668
$ws.Cells.Item(13,15).Value = 0.2224806037672029
$ws.Cells.Item(13,16).Value = 0.2224806037672029
$ws.Cells.Item(13,17).Value = 5.887477803199999
$ws.Cells.Item(13,18).Value = 52.9873002288
$ws.Cells.Item(13,19).Value = 0.05748570030832437
$ws.Cells.Item(13,20).Value = 0.05748570030832437

# Row 14
$ws.Cells.Item(14,1).Value = "sCs"
$ws.Cells.Item(14,2).Value = "Dnajb11"
$ws.Cells.Item(14,3).Value = "Prtg"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 5.953784333333334
$ws.Cells.Item(14,8).Value = 17.861353
$ws.Cells.Item(14,9).Value = 0.06416365763943291
$ws.Cells.Item(14,10).Value = 0.0641636576394329
$ws.Cells.Item(14,11).Value = 2
$ws.Cells.Item(14,12).Value = 0.6666666666666666
$ws.Cells.Item(14,13).Value = 0.130725
$ws.Cells.Item(14,14).Value = 0.392175
$ws.Cells.Item(14,15).Value = 0.1184385768344502
$ws.Cells.Item(14,16).Value = 0.1184385768344502
$ws.Cells.Item(14,17).Value = 0.7783084569750002
$ws.Cells.Item(14,18).Value = 7.004776112775001
$ws.Cells.Item(14,19).Value = 0.007599452295307335
$ws.Cells.Item(14,20).Value = 0.007599452295307334

# Row 15
$ws.Cells.Item(15,1).Value = "sCs"
$ws.Cells.Item(15,2).Value = "Dnajb11"
$ws.Cells.Item(15,3).Value = "Prtg"
$ws.Cells.Item(15,4).Value = "FAPs"
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 5.953784333333334
$ws.Cells.Item(15,8).Value = 17.861353
$ws.Cells.Item(15,9).Value = 0.06416365763943291
$ws.Cells.Item(15,10).Value = 0.0641636576394329
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 0.7274516666666667
$ws.Cells.Item(15,14).Value = 2.182355
$ws.Cells.Item(15,15).Value = 0.6590808193983468
$ws.Cells.Item(15,16).Value = 0.659080819398347
$ws.Cells.Item(15,17).Value = 4.331090336257223
$ws.Cells.Item(15,18).Value = 38.97981302631501
$ws.Cells.Item(15,19).Value = 0.04228903605259244
$ws.Cells.Item(15,20).Value = 0.04228903605259244

# Row 16
$ws.Cells.Item(16,1).Value = "sCs"
$ws.Cells.Item(16,2).Value = "Dnajb11"
$ws.Cells.Item(16,3).Value = "Prtg"
$ws.Cells.Item(16,4).Value = "sCs"
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 5.953784333333334
$ws.Cells.Item(16,8).Value = 17.861353
$ws.Cells.Item(16,9).Value = 0.06416365763943291
$ws.Cells.Item(16,10).Value = 0.0641636576394329
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 0.24556
$ws.Cells.Item(16,14).Value = 0.73668
$ws.Cells.Item(16,15).Value = 0.2224806037672029
$ws.Cells.Item(16,16).Value = 0.2224806037672029
$ws.Cells.Item(16,17).Value = 1.462011280893333
$ws.Cells.Item(16,18).Value = 13.15810152804
$ws.Cells.Item(16,19).Value = 0.01427516929153313
$ws.Cells.Item(16,20).Value = 0.01427516929153313

